# Apply label flips recorded in the source diff (column A, single-column label sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each tuple is (row, newValue) for cell A<row>.
$changes = @(
    @(214, 0),
    @(236, 1),
    @(249, 0),
    @(298, 1),
    @(316, 0),
    @(330, 1),
    @(353, 1),
    @(362, 0),
    @(387, 0),
    @(433, 1),
    @(442, 1),
    @(721, 0),
    @(722, 0),
    @(723, 0),
    @(724, 0),
    @(725, 0),
    @(727, 0),
    @(728, 0),
    @(729, 0),
    @(730, 0),
    @(733, 0),
    @(737, 0),
    @(740, 0),
    @(743, 0),
    @(746, 0),
    @(751, 1),
    @(759, 1),
    @(761, 1),
    @(767, 0),
    @(768, 1),
    @(769, 1),
    @(775, 1),
    @(777, 1),
    @(778, 0),
    @(782, 0),
    @(784, 1),
    @(786, 1),
    @(791, 0),
    @(793, 0),
    @(797, 1),
    @(803, 0),
    @(806, 1),
    @(810, 1),
    @(819, 0),
    @(821, 0),
    @(823, 0),
    @(824, 0),
    @(827, 1),
    @(829, 1),
    @(837, 1),
    @(839, 0),
    @(845, 1),
    @(849, 0),
    @(854, 0),
    @(855, 1),
    @(857, 0),
    @(861, 0),
    @(864, 0),
    @(865, 1),
    @(878, 0),
    @(879, 1),
    @(881, 1),
    @(882, 1),
    @(885, 1),
    @(892, 1),
    @(893, 1),
    @(894, 1),
    @(895, 1),
    @(897, 1),
    @(898, 1),
    @(899, 1),
    @(900, 1),
    @(901, 1),
    @(902, 1),
    @(906, 1),
    @(907, 1),
    @(1050, 0),
    @(1082, 0),
    @(1219, 0),
    @(1288, 1),
    @(1301, 0),
    @(1328, 1),
    @(1363, 0),
    @(1386, 1),
    @(1424, 1),
    @(1633, 0),
    @(1636, 1),
    @(1642, 0),
    @(1645, 1),
    @(1657, 0),
    @(1659, 0),
    @(1662, 0),
    @(1667, 1),
    @(1668, 0),
    @(1669, 0),
    @(1672, 0),
    @(1678, 1),
    @(1705, 0),
    @(1707, 1),
    @(1730, 0),
    @(1736, 1),
    @(1739, 0),
    @(1750, 1),
    @(1756, 1),
    @(1759, 1),
    @(1762, 0),
    @(1766, 0),
    @(1772, 1),
    @(1778, 0),
    @(1780, 0),
    @(1782, 1),
    @(1794, 1),
    @(1797, 0),
    @(1798, 1),
    @(1799, 0),
    @(1800, 0)
)

foreach ($pair in $changes) {
    $row = $pair[0]
    $newVal = $pair[1]
    $ws.Cells.Item($row, 1).Value = $newVal
}

Write-Host "Applied $($changes.Count) cell updates to column A"
